$d = $word.ActiveDocument

# --- locate the target paragraph: "O Projeto será uma árvore de decisões..." ---
$target = $null
foreach ($p in $d.Paragraphs) {
    if ($p.Range.Text.StartsWith("O Projeto")) {
        $target = $p
    }
}

$pStart = $target.Range.Start
$pEnd = $target.Range.End

# The paragraph currently holds two runs: "O" and " Projeto será uma árvore de
# decisões para fazer as classificações de alunos. ". Clear everything after
# the leading "O" so we can rebuild the trailing runs to match the new text.
$clearRange = $d.Range($pStart + 1, $pEnd - 1)
$clearRange.Text = ""

# --- rebuild the trailing runs, one InsertAfter per run, each stamped with
#     the same font size (12pt == <w:sz w:val="24"/>) as the rest of the
#     paragraph so a genuine new <w:r> is minted instead of being merged
#     back into a neighbour. ---
function Append-Run([string]$text) {
    $endPos = $target.Range.End - 1
    $anchor = $d.Range($endPos, $endPos)
    $anchor.InsertAfter($text)
    $newRange = $d.Range($endPos, $endPos + $text.Length)
    $newRange.Font.Size = 12
}

Append-Run(" Projeto era")
Append-Run(" uma árvore de decisões para fazer as classificações de alunos.")
Append-Run(" Mas devido ao tema não ser muito adaptável a este tipo de árvores, o projeto será feito com uma árvore binária.")

# --- move the _GoBack bookmark from the "Descrição:" heading paragraph down
#     to the end of this paragraph (Bookmarks.Add replaces any existing
#     bookmark of the same name, so the old one disappears automatically). ---
$bmEnd = $target.Range.End - 1
$bmRange = $d.Range($bmEnd - 1, $bmEnd)
$d.Bookmarks.Add("_GoBack", $bmRange)

Write-Output "done"
